$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-04-15 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-16 Tuesday", 2) | Out-Null

# Update each answer cell in the practice table (row-major order, 20 rows x 5 cols)
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "18+63=81"
$t.Rows.Item(1).Cells.Item(2).Range.Text = "66-30=36"
$t.Rows.Item(1).Cells.Item(3).Range.Text = "42+3=45"
$t.Rows.Item(1).Cells.Item(4).Range.Text = "66+5=71"
$t.Rows.Item(1).Cells.Item(5).Range.Text = "10+54=64"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "47-13=34"
$t.Rows.Item(2).Cells.Item(2).Range.Text = "66-11=55"
$t.Rows.Item(2).Cells.Item(3).Range.Text = "88-3=85"
$t.Rows.Item(2).Cells.Item(4).Range.Text = "79-30=49"
$t.Rows.Item(2).Cells.Item(5).Range.Text = "19-5=14"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "99-67=32"
$t.Rows.Item(3).Cells.Item(2).Range.Text = "88-34=54"
$t.Rows.Item(3).Cells.Item(3).Range.Text = "70-35=35"
$t.Rows.Item(3).Cells.Item(4).Range.Text = "93+5=98"
$t.Rows.Item(3).Cells.Item(5).Range.Text = "71+11=82"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "76-61=15"
$t.Rows.Item(4).Cells.Item(2).Range.Text = "17+43=60"
$t.Rows.Item(4).Cells.Item(3).Range.Text = "42+53=95"
$t.Rows.Item(4).Cells.Item(4).Range.Text = "31+30=61"
$t.Rows.Item(4).Cells.Item(5).Range.Text = "61-35=26"
$t.Rows.Item(5).Cells.Item(1).Range.Text = "39-23=16"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "56+34=90"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "50+3=53"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "78+2=80"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "7+66=73"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "28+66=94"
$t.Rows.Item(6).Cells.Item(2).Range.Text = "15-13=2"
$t.Rows.Item(6).Cells.Item(3).Range.Text = "58-13=45"
$t.Rows.Item(6).Cells.Item(4).Range.Text = "95-24=71"
$t.Rows.Item(6).Cells.Item(5).Range.Text = "97+1=98"
$t.Rows.Item(7).Cells.Item(1).Range.Text = "75-63=12"
$t.Rows.Item(7).Cells.Item(2).Range.Text = "54+3=57"
$t.Rows.Item(7).Cells.Item(3).Range.Text = "69-45=24"
$t.Rows.Item(7).Cells.Item(4).Range.Text = "85-56=29"
$t.Rows.Item(7).Cells.Item(5).Range.Text = "27+24=51"
$t.Rows.Item(8).Cells.Item(1).Range.Text = "81-0=81"
$t.Rows.Item(8).Cells.Item(2).Range.Text = "7+39=46"
$t.Rows.Item(8).Cells.Item(3).Range.Text = "37+16=53"
$t.Rows.Item(8).Cells.Item(4).Range.Text = "53+39=92"
$t.Rows.Item(8).Cells.Item(5).Range.Text = "49-23=26"
$t.Rows.Item(9).Cells.Item(1).Range.Text = "53-20=33"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "15+29=44"
$t.Rows.Item(9).Cells.Item(3).Range.Text = "17-4=13"
$t.Rows.Item(9).Cells.Item(4).Range.Text = "32-0=32"
$t.Rows.Item(9).Cells.Item(5).Range.Text = "29-27=2"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "48-0=48"
$t.Rows.Item(10).Cells.Item(2).Range.Text = "99-37=62"
$t.Rows.Item(10).Cells.Item(3).Range.Text = "41-12=29"
$t.Rows.Item(10).Cells.Item(4).Range.Text = "11+46=57"
$t.Rows.Item(10).Cells.Item(5).Range.Text = "42+10=52"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "24+69=93"
$t.Rows.Item(11).Cells.Item(2).Range.Text = "93-24=69"
$t.Rows.Item(11).Cells.Item(3).Range.Text = "96-17=79"
$t.Rows.Item(11).Cells.Item(4).Range.Text = "34-27=7"
$t.Rows.Item(11).Cells.Item(5).Range.Text = "12+22=34"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "95-28=67"
$t.Rows.Item(12).Cells.Item(2).Range.Text = "70-53=17"
$t.Rows.Item(12).Cells.Item(3).Range.Text = "84-28=56"
$t.Rows.Item(12).Cells.Item(4).Range.Text = "62-27=35"
$t.Rows.Item(12).Cells.Item(5).Range.Text = "14+33=47"
$t.Rows.Item(13).Cells.Item(1).Range.Text = "97-12=85"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "3+26=29"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "14+46=60"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "84+8=92"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "56-34=22"
$t.Rows.Item(14).Cells.Item(1).Range.Text = "3+86=89"
$t.Rows.Item(14).Cells.Item(2).Range.Text = "60+0=60"
$t.Rows.Item(14).Cells.Item(3).Range.Text = "2+22=24"
$t.Rows.Item(14).Cells.Item(4).Range.Text = "71+26=97"
$t.Rows.Item(14).Cells.Item(5).Range.Text = "37-6=31"
$t.Rows.Item(15).Cells.Item(1).Range.Text = "46+6=52"
$t.Rows.Item(15).Cells.Item(2).Range.Text = "26+32=58"
$t.Rows.Item(15).Cells.Item(3).Range.Text = "5+71=76"
$t.Rows.Item(15).Cells.Item(4).Range.Text = "81-30=51"
$t.Rows.Item(15).Cells.Item(5).Range.Text = "28-9=19"
$t.Rows.Item(16).Cells.Item(1).Range.Text = "60-6=54"
$t.Rows.Item(16).Cells.Item(2).Range.Text = "77-40=37"
$t.Rows.Item(16).Cells.Item(3).Range.Text = "42-32=10"
$t.Rows.Item(16).Cells.Item(4).Range.Text = "66-15=51"
$t.Rows.Item(16).Cells.Item(5).Range.Text = "76-8=68"
$t.Rows.Item(17).Cells.Item(1).Range.Text = "50+37=87"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "94-32=62"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "35+43=78"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "76-4=72"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "81-10=71"
$t.Rows.Item(18).Cells.Item(1).Range.Text = "48-23=25"
$t.Rows.Item(18).Cells.Item(2).Range.Text = "58-15=43"
$t.Rows.Item(18).Cells.Item(3).Range.Text = "40-31=9"
$t.Rows.Item(18).Cells.Item(4).Range.Text = "92-36=56"
$t.Rows.Item(18).Cells.Item(5).Range.Text = "65+20=85"
$t.Rows.Item(19).Cells.Item(1).Range.Text = "61-34=27"
$t.Rows.Item(19).Cells.Item(2).Range.Text = "65-51=14"
$t.Rows.Item(19).Cells.Item(3).Range.Text = "93-15=78"
$t.Rows.Item(19).Cells.Item(4).Range.Text = "11+48=59"
$t.Rows.Item(19).Cells.Item(5).Range.Text = "32+63=95"
$t.Rows.Item(20).Cells.Item(1).Range.Text = "34-30=4"
$t.Rows.Item(20).Cells.Item(2).Range.Text = "63-25=38"
$t.Rows.Item(20).Cells.Item(3).Range.Text = "51+41=92"
$t.Rows.Item(20).Cells.Item(4).Range.Text = "42+50=92"
$t.Rows.Item(20).Cells.Item(5).Range.Text = "35+0=35"
